$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the hyperlinks anchored on A3 and B3 only (leave A2/B2 hyperlinks intact).
# Deleting via Range("A3:B3").Hyperlinks.Delete() removes every hyperlink on the
# sheet, so instead walk the sheet's Hyperlinks collection and delete just the
# matching ones, re-scanning after each delete since the collection re-indexes.
$changed = $true
while ($changed) {
    $changed = $false
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq '$A$3' -or $addr -eq '$B$3') {
            $hl.Delete()
            $changed = $true
            break
        }
    }
}

# Clear the now un-linked cell values, but keep their (Hyperlink) style.
$ws.Range("A3:B3").ClearContents()

# Update the active selection (was D5) to A3.
$ws.Range("A3").Select()
